$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1587.375
$ws.Range("J17").Value = 1587.375
$ws.Range("L17").Value = 4762.125
$ws.Range("N17").Value = -5098.125
$ws.Range("H58").Value = 22117.34
$ws.Range("J58").Value = 26258.5
$ws.Range("L58").Value = 78775.5
$ws.Range("N58").Value = -79075.5
$ws.Range("H132").Value = 4257789
$ws.Range("I132").Value = 4880199.5
$ws.Range("K132").Value = 14640598.5
$ws.Range("M132").Value = -14638068.5
$ws.Range("H135").Value = 669
$ws.Range("I135").Value = 488.42554
$ws.Range("J135").Value = 1612
$ws.Range("K135").Value = 4395.82986
$ws.Range("L135").Value = 14508
$ws.Range("M135").Value = -1860.82986
$ws.Range("N135").Value = -19578
$ws.Range("H136").Value = 31035
$ws.Range("J136").Value = 31035
$ws.Range("L136").Value = 31035
$ws.Range("N136").Value = -41235
$ws.Range("H137").Value = 2626.739
$ws.Range("I137").Value = 2621.2
$ws.Range("J137").Value = 2644.3635
$ws.Range("K137").Value = 7863.599999999999
$ws.Range("L137").Value = 7933.0905
$ws.Range("M137").Value = -5313.599999999999
$ws.Range("N137").Value = -13033.0905
$ws.Range("H138").Value = 1792.4
$ws.Range("I138").Value = 531.12195
$ws.Range("J138").Value = 2668.8813
$ws.Range("K138").Value = 1593.36585
$ws.Range("L138").Value = 8006.6439
$ws.Range("M138").Value = 3546.63415
$ws.Range("N138").Value = -18286.6439

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6375.677
$ws.Range("I32").Value = 5189.579
$ws.Range("J32").Value = 10294.956
$ws.Range("K32").Value = 5189.579
$ws.Range("L32").Value = 10294.956
$ws.Range("M32").Value = -4902.579
$ws.Range("N32").Value = -10868.956
$ws.Range("H61").Value = 2504.9285
$ws.Range("I61").Value = 796.7368
$ws.Range("K61").Value = 796.7368
$ws.Range("M61").Value = -584.7368
$ws.Range("H74").Value = 483.13727
$ws.Range("I74").Value = 440.41666
$ws.Range("J74").Value = 1166.6666
$ws.Range("K74").Value = 440.41666
$ws.Range("L74").Value = 1166.6666
$ws.Range("M74").Value = 433.58334
$ws.Range("N74").Value = -2914.6666
$ws.Range("H77").Value = 483.13727
$ws.Range("I77").Value = 440.41666
$ws.Range("J77").Value = 1166.6666
$ws.Range("K77").Value = 2202.0833
$ws.Range("L77").Value = 5833.333000000001
$ws.Range("M77").Value = 2165.9167
$ws.Range("N77").Value = -14569.333
$ws.Range("H136").Value = 2504.9285
$ws.Range("I136").Value = 796.7368
$ws.Range("K136").Value = 2390.2104
$ws.Range("M136").Value = 159.7896000000001
$ws.Range("H138").Value = 98966.664
$ws.Range("J138").Value = 98966.664
$ws.Range("L138").Value = 98966.664
$ws.Range("N138").Value = -109246.664

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2082.7273
$ws.Range("I31").Value = 1334.2979
$ws.Range("J31").Value = 3255.2666
$ws.Range("K31").Value = 1334.2979
$ws.Range("L31").Value = 3255.2666
$ws.Range("M31").Value = -1039.2979
$ws.Range("N31").Value = -3845.2666
$ws.Range("H34").Value = 2082.7273
$ws.Range("I34").Value = 1334.2979
$ws.Range("J34").Value = 3255.2666
$ws.Range("K34").Value = 1334.2979
$ws.Range("L34").Value = 3255.2666
$ws.Range("M34").Value = -1132.2979
$ws.Range("N34").Value = -3659.2666
$ws.Range("H58").Value = 7354849
$ws.Range("I58").Value = 904.2453
$ws.Range("J58").Value = 33338788
$ws.Range("K58").Value = 904.2453
$ws.Range("L58").Value = 33338788
$ws.Range("M58").Value = -701.2453
$ws.Range("N58").Value = -33339194
$ws.Range("H99").Value = 2708.0908
$ws.Range("I99").Value = 1786.25
$ws.Range("K99").Value = 1786.25
$ws.Range("M99").Value = -288.25
$ws.Range("H107").Value = 1957.6875
$ws.Range("I107").Value = 1147.5454
$ws.Range("J107").Value = 3740
$ws.Range("K107").Value = 1147.5454
$ws.Range("L107").Value = 3740
$ws.Range("M107").Value = 772.4546
$ws.Range("N107").Value = -7580
$ws.Range("H126").Value = 2708.0908
$ws.Range("I126").Value = 1786.25
$ws.Range("K126").Value = 5358.75
$ws.Range("M126").Value = -2888.75
$ws.Range("H132").Value = 1605.0741
$ws.Range("I132").Value = 1199.3077
$ws.Range("J132").Value = 2660.0667
$ws.Range("K132").Value = 3597.9231
$ws.Range("L132").Value = 7980.2001
$ws.Range("M132").Value = -1067.9231
$ws.Range("N132").Value = -13040.2001
$ws.Range("H134").Value = 1760.9286
$ws.Range("I134").Value = 699.65216
$ws.Range("K134").Value = 2098.95648
$ws.Range("M134").Value = 436.0435200000002
$ws.Range("H136").Value = 7354849
$ws.Range("I136").Value = 904.2453
$ws.Range("J136").Value = 33338788
$ws.Range("K136").Value = 2712.7359
$ws.Range("L136").Value = 100016364
$ws.Range("M136").Value = -162.7359000000001
$ws.Range("N136").Value = -100021464

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1218.091
$ws.Range("I97").Value = 593
$ws.Range("J97").Value = 1452.5
$ws.Range("K97").Value = 1779
$ws.Range("L97").Value = 4357.5
$ws.Range("M97").Value = -1283
$ws.Range("N97").Value = -5349.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 94044.63
$ws.Range("I102").Value = 2060.4
$ws.Range("J102").Value = 170698.17
$ws.Range("K102").Value = 2060.4
$ws.Range("L102").Value = 170698.17
$ws.Range("M102").Value = -438.4000000000001
$ws.Range("N102").Value = -173942.17
$ws.Range("H132").Value = 2243.86
$ws.Range("I132").Value = 1738.1428
$ws.Range("J132").Value = 2887.5
$ws.Range("K132").Value = 5214.428400000001
$ws.Range("L132").Value = 8662.5
$ws.Range("M132").Value = -2684.428400000001
$ws.Range("N132").Value = -13722.5
$ws.Range("H136").Value = 16452.334
$ws.Range("J136").Value = 16452.334
$ws.Range("L136").Value = 49357.00199999999
$ws.Range("N136").Value = -54457.00199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2555.6667
$ws.Range("I46").Value = 600.3333
$ws.Range("J46").Value = 3533.3333
$ws.Range("K46").Value = 600.3333
$ws.Range("L46").Value = 3533.3333
$ws.Range("M46").Value = -412.3333
$ws.Range("N46").Value = -3909.3333
$ws.Range("H132").Value = 1694.6666
$ws.Range("I132").Value = 995.5111000000001
$ws.Range("K132").Value = 2986.5333
$ws.Range("M132").Value = -456.5333000000001
$ws.Range("H135").Value = 29577.75
$ws.Range("J135").Value = 29577.75
$ws.Range("L135").Value = 29577.75
$ws.Range("N135").Value = -39717.75
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 30000
$ws.Range("J46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30462
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -95070
$ws.Range("H140").Value = 49999.4
$ws.Range("J140").Value = 49999.4
$ws.Range("L140").Value = 49999.4
$ws.Range("N140").Value = -60359.4
